# Regenerate the localization-status report for the new handoff file
# (69f57c00-3142-43d1-8e83-b210ad9ab90c -> 11344c0d-8125-4cf9-bcd1-829975378b72)
# and refresh the "Ready for handoff" timestamps / xliff handback info.

$wb = $excel.ActiveWorkbook

$oldGuid = "69f57c00-3142-43d1-8e83-b210ad9ab90c"
$newGuid = "11344c0d-8125-4cf9-bcd1-829975378b72"
$newHash = "672c0e9881d01bd882a5cbbdd7742359378dce0d"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid.md"

$b2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b498c8a826b3da58c039a7f17d0c3a35c8d7cb9a/e2e/$oldGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
$ws.Hyperlinks.Add($ws.Range("B2"), $b2Addr, "", "", "e2e\$newGuid.md")

$ws.Range("G2").Value = "2016-09-02 13:09:42"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$a2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b498c8a826b3da58c039a7f17d0c3a35c8d7cb9a/e2e/$oldGuid.md"
$ws.Range("A2").Value = "$newGuid.md"
$ws.Hyperlinks.Add($ws.Range("A2"), $a2Addr, "", "", "$newGuid.md")

$ws.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.75
$ws.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$a2Addr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b498c8a826b3da58c039a7f17d0c3a35c8d7cb9a/e2e/$oldGuid.md"
$ws.Range("A2").Value = "$newGuid.md"
$ws.Hyperlinks.Add($ws.Range("A2"), $a2Addr, "", "", "$newGuid.md")

$ws.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-09-02 13:09:42"

$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""

$ws.Range("K2").Value = "0001-01-01 00:00:00"

$ws.Columns.Item(9).ColumnWidth = 17.75
$ws.Columns.Item(10).ColumnWidth = 20.75
